# Update the cryptos list data (prices / hourly volume changes) to match
# the latest scrape. Numeric-looking price strings are prefixed with a
# leading apostrophe so Excel keeps storing them as text (matching the
# source data, which stores every Price cell as text) instead of
# auto-converting them into floating point numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.963.35"
$ws.Range("E2").Value = "  +0.12%  "
$ws.Range("D3").Value = "1.560.77"
$ws.Range("E3").Value = "  +0.45%  "
$ws.Range("D5").Value = "'207.27"
$ws.Range("E5").Value = "  +0.10%  "
$ws.Range("E6").Value = "  +0.39%  "
$ws.Range("D8").Value = "'22.11"
$ws.Range("E8").Value = "  +1.59%  "
$ws.Range("E9").Value = "  +0.07%  "
$ws.Range("D10").Value = "'0.0599"
$ws.Range("E10").Value = "  +2.15%  "
$ws.Range("E11").Value = "  +0.02%  "
$ws.Range("D12").Value = "1.782.94"
$ws.Range("E12").Value = "  +0.43%  "
$ws.Range("D13").Value = "1.561.84"
$ws.Range("E13").Value = "  +0.43%  "
$ws.Range("E14").Value = "  +0.59%  "
$ws.Range("D16").Value = "'62.13"
$ws.Range("E16").Value = "  +0.70%  "
$ws.Range("D17").Value = "26.972.31"
$ws.Range("E17").Value = "  +0.15%  "
$ws.Range("D18").Value = "'217.16"
$ws.Range("E18").Value = "  +0.04%  "
$ws.Range("D19").Value = "0.0₃0704"
$ws.Range("E19").Value = "  +2.33%  "
$ws.Range("E20").Value = "  +1.77%  "
$ws.Range("E21").Value = "  -0.22%  "
$ws.Range("E22").Value = "  +1.67%  "
$ws.Range("D23").Value = "'9.20"
$ws.Range("E23").Value = "  -0.02%  "
$ws.Range("E24").Value = "  -1.35%  "
$ws.Range("D25").Value = "'153.39"
$ws.Range("E25").Value = "  -0.32%  "
$ws.Range("D26").Value = "'6.62"
$ws.Range("E26").Value = "  +0.68%  "
$ws.Range("D27").Value = "'15.09"
$ws.Range("E27").Value = "  +1.55%  "
$ws.Range("D28").Value = "'0.105"
$ws.Range("E28").Value = "  +1.52%  "
$ws.Range("E29").Value = "  -0.39%  "
$ws.Range("E30").Value = "  +0.36%  "
$ws.Range("E31").Value = "  +1.36%  "
$ws.Range("E33").Value = "  +3.34%  "
$ws.Range("D34").Value = "1.421.19"
$ws.Range("E34").Value = "  -0.29%  "
$ws.Range("E35").Value = "  +3.29%  "
$ws.Range("D36").Value = "'1.05"
$ws.Range("E36").Value = "  +9.21%  "
$ws.Range("E37").Value = "  +1.52%  "
$ws.Range("E38").Value = "  +0.70%  "
$ws.Range("D39").Value = "'0.532"
$ws.Range("E39").Value = "  +2.40%  "
$ws.Range("E40").Value = "  -0.09%  "
$ws.Range("E41").Value = "  -0.23%  "
$ws.Range("E42").Value = "  -0.30%  "
$ws.Range("B43").Value = "WEMIXToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D43").Value = "'1.01"
$ws.Range("E43").Value = "  +1.98%  "
$ws.Range("B44").Value = "MXToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D44").Value = "'2.33"
$ws.Range("E44").Value = "  +2.43%  "
$ws.Range("D45").Value = "'64.77"
$ws.Range("E45").Value = "  +1.54%  "
$ws.Range("D46").Value = "'1.75"
$ws.Range("E46").Value = "  +0.60%  "
$ws.Range("D47").Value = "1.696.82"
$ws.Range("E47").Value = "  +0.41%  "
$ws.Range("D48").Value = "'87.49"
$ws.Range("E48").Value = "  +1.56%  "
$ws.Range("E49").Value = "  -0.34%  "
$ws.Range("B50").Value = "BabyDogeCoin"
$ws.Range("C50").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D50").Value = "0.0₆0101"
$ws.Range("E50").Value = "  +0.78%  "
$ws.Range("B51").Value = "Algorand"
$ws.Range("C51").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D51").Value = "'0.0958"
$ws.Range("E51").Value = "  +0.17%  "
